$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, shifting rows 4-101 down to 5-102
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with its new data
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44812
$ws.Cells.Item(4, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100102
$ws.Cells.Item(4, 8).Value = "Cítricos"
$ws.Cells.Item(4, 9).Value = 100102005
$ws.Cells.Item(4, 10).Value = "Naranja"
$ws.Cells.Item(4, 11).Value = "Fukumoto"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 300
$ws.Cells.Item(4, 14).Value = 600
$ws.Cells.Item(4, 15).Value = 650
$ws.Cells.Item(4, 16).Value = 625
$ws.Cells.Item(4, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(4, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(4, 19).Value = 625
$ws.Cells.Item(4, 20).Value = 1
